$wb = $excel.ActiveWorkbook

# Mapping of old sheet names -> new sheet names (in workbook order / rId order)
$newNames = @(
    "summ13356575",
    "summ13615682",
    "summ13865359",
    "summ14131831",
    "summ14387273",
    "summ14691676",
    "summ14954215",
    "summ15218137",
    "summ15489248"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename the sheet
    $ws.Name = $newNames[$i - 1]

    # Update the Education[T.Unknown] label to Education[T.Unknown/Other] in cell A5
    if ($ws.Range("A5").Text -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
